$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'50.036.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.49%  "
$ws.Range("D3").Value = "'2.667.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +7.62%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'114.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.86%  "
$ws.Range("D6").Value = "'326.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.11%  "
$ws.Range("D7").Value = "'0.528"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.12%  "
$ws.Range("D9").Value = "'0.557"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.40%  "
$ws.Range("D10").Value = "'41.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.74%  "
$ws.Range("D11").Value = "'20.10"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "'0.0825"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.45%  "
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").Value = "'7.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.98%  "
$ws.Range("D15").Value = "'3.086.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.57%  "
$ws.Range("D16").Value = "'2.681.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.82%  "
$ws.Range("D17").Value = "'0.877"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.51%  "
$ws.Range("D18").Value = "'49.960.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.47%  "
$ws.Range("D19").Value = "'13.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.94%  "
$ws.Range("D20").Value = "'6.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.47%  "
$ws.Range("D21").Value = "'2.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").Value = "'0.0₃0961"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.98%  "
$ws.Range("D23").Value = "'72.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.83%  "
$ws.Range("D24").Value = "'277.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").Value = "'2.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.38%  "
$ws.Range("D26").Value = "'26.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.31%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "'10.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.23%  "
$ws.Range("D29").Value = "'2.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").Value = "'36.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.71%  "
$ws.Range("E31").Value = "  +4.66%  "
$ws.Range("D32").Value = "'50.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("D33").Value = "'5.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.74%  "
$ws.Range("D34").Value = "'19.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.17%  "
$ws.Range("D35").Value = "'0.0819"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.65%  "
$ws.Range("D36").Value = "'5.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.82%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  +7.72%  "
$ws.Range("E39").Value = "  +10.04%  "
$ws.Range("D40").Value = "'125.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.06%  "
$ws.Range("D41").Value = "'0.113"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.55%  "
$ws.Range("D42").Value = "'22.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.31%  "
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("D44").Value = "'0.0319"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.45%  "
$ws.Range("D45").Value = "'2.115.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.22%  "
$ws.Range("D46").Value = "'3.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.70%  "
$ws.Range("E47").Value = "  +14.96%  "
$ws.Range("D48").Value = "'2.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.39%  "
$ws.Range("D49").Value = "'9.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.72%  "
$ws.Range("D50").Value = "'5.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.34%  "
$ws.Range("D51").Value = "'59.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.98%  "
